$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet.
$ws.Name = "TC001_TestCase_001"

# The sheet had a duplicate "ENV" column (C) mirroring column B ("ENV1").
# Drop column C entirely so column B (renamed "ENV") is the only env column.
$ws.Columns("C:C").Delete()
$ws.Range("B1").Value = "ENV"

# Move the active selection as recorded in the saved sheet view.
$ws.Range("B35").Select() | Out-Null
